$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 50, shifting existing rows 50-140 down to 51-141.
$ws.Rows("50:50").Insert()

# Populate the new row 50 with the new weekly data point.
$ws.Range("A50").Value = 4
$ws.Range("B50").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C50").Value = "Los Lagos"
$ws.Range("D50").Value = 44469
$ws.Range("E50").Value = 10
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100108
$ws.Range("H50").Value = "Tropicales y subtropicales"
$ws.Range("I50").Value = 100108005
$ws.Range("J50").Value = "Piña"
$ws.Range("K50").Value = "Caramelo"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 80
$ws.Range("N50").Value = 22000
$ws.Range("O50").Value = 22000
$ws.Range("P50").Value = 22000
$ws.Range("Q50").Value = "`$/caja 12 unidades"
$ws.Range("R50").Value = "Ecuador"
$ws.Range("S50").Value = 1833
$ws.Range("T50").Value = 12
